# Update "Forecast Comparison" sheet with correct forecast output:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - shorten the "Week" labels (W01 -> W1, etc.)
#  - store the is_holiday_week column as boolean values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column before column B ("ASIN").
#    Everything currently in B:I shifts right to C:J.
$ws.Columns.Item(2).Insert()

# 2. New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# 3. Fill in the Week_Start_Date values (column B, rows 2-17).
#    Force text formatting first so Excel doesn't coerce the
#    "yyyy-mm-dd" strings into date serial numbers, then clear
#    the temporary formatting so the cells keep the default style.
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

$dateRange = $ws.Range("B2:B17")
$dateRange.NumberFormat = "@"
for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $weekStartDates[$i]
}
$dateRange.ClearFormats()

# 4. Shorten the Week labels in column A (W01 -> W1 ... W09 -> W9).
#    W10-W16 are already in the right shape.
$weekLabels = @("W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8", "W9")
for ($i = 0; $i -lt $weekLabels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $weekLabels[$i]
}

# 5. is_holiday_week (now column J after the insert) becomes boolean.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}
